$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows correspond to worksheet rows 2-25 (index column A = 0..23)
# Columns updated: B, D, E, F, G, K
$data = @(
    @(7.826708101776549, 8.428653024961088, 16.31824107584387, 46.98426329439417, 3.671752231996937, 13.56665419506214),
    @(7.758101565779985, 8.280758455761747, 15.39753157995865, 45.63964592489109, 3.676753955401931, 13.28718078769574),
    @(7.717641409760542, 8.188755721533065, 14.80962437884951, 44.80113486122085, 3.67997244360572, 13.12080524575059),
    @(7.701589100054782, 8.150989038213266, 14.56464244413985, 44.45662460226387, 3.681321265927033, 13.05444171113893),
    @(7.698950410969188, 8.144702069222351, 14.52364601513326, 44.39926246137694, 3.681547493085362, 13.04351234774442),
    @(7.717423137423782, 8.188247465365814, 14.80634194923813, 44.79649945829365, 3.679990483162585, 13.11990427009909),
    @(7.802717522349337, 8.377922182929384, 16.0056177458328, 46.52357234892643, 3.673446349161892, 13.46928219456833),
    @(7.982404126339218, 8.739300674617681, 18.19951367570059, 49.79014626091607, 3.661774022430802, 14.19037544236281),
    @(8.12090192165136, 8.996955917567968, 19.84232061102844, 52.09508635747896, 3.653893390247799, 14.73479051535524),
    @(8.18507392196563, 9.112186251827218, 20.54901404399155, 53.11903725904082, 3.650456483817698, 14.984140353378),
    @(8.209521585795459, 9.155514392211323, 20.81082156078195, 53.50295315353117, 3.649176089380672, 15.07868783092537),
    @(8.204250134732963, 9.14619690723336, 20.75469403046332, 53.42044486237781, 3.649450910550092, 15.05832165608762),
    @(8.187082417615159, 9.115757163797387, 20.57066915507828, 53.15070065383571, 3.650350723520893, 14.99191700613855),
    @(8.17658523205591, 9.097071299332619, 20.45719421485452, 52.9849671453202, 3.650904625210135, 14.9512549634175),
    @(8.116729859724844, 8.98938357647334, 19.79532244295801, 52.02764806088304, 3.654120961937525, 14.71851864715051),
    @(8.080294990255656, 8.922797501326727, 19.37890848544015, 51.43384759569712, 3.656131853512537, 14.57608929453312),
    @(8.059449879577045, 8.884314023403777, 19.13557171560259, 51.09001179823385, 3.657302410701852, 14.49433382162876),
    @(8.052411828850287, 8.871253125073144, 19.05252363906764, 50.97320986660721, 3.657701142466935, 14.46668508472755),
    @(8.084162180407484, 8.929905009525923, 19.42363219583508, 51.49729869904793, 3.655916348802328, 14.59123483196193),
    @(8.192121166644206, 9.124706562923665, 20.62487886438347, 53.23003729742253, 3.650085855928148, 15.01141921448411),
    @(8.263526402034294, 9.250221055578686, 21.37616347181789, 54.34002797354653, 3.646398111994985, 15.2867064138488),
    @(8.225345437477033, 9.183403434341521, 20.97826810256213, 53.74975052946004, 3.648355158520299, 15.13975676497951),
    @(8.082413505128102, 8.926692333324509, 19.40342486555761, 51.46862006974363, 3.656013733399362, 14.58438712338009),
    @(7.932574751704799, 8.642818370300228, 17.60480309248785, 48.92167977231189, 3.664808730723989, 13.99221123460496)
)

$columns = @("B", "D", "E", "F", "G", "K")

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $cellRef = "" + $columns[$j] + $rowNum
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
